$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 635, shifting the existing rows
# (635..676) down to (636..677). This matches the diff, which
# shows a new "2026/01/16" / "金" / 14 / 58 record spliced in just
# before the "2026/12/29" block, with everything after it unchanged
# but renumbered one row lower.
$ws.Rows(635).Insert()

# Column A stores dates as literal text (e.g. "2026/01/16"), not real
# date serials. Pre-format the new cell as Text so Excel doesn't
# auto-convert the string into a date value/format, then clear the
# formatting back off so the cell ends up styled like its neighbours
# (no explicit style) while keeping the literal text value.
$ws.Range("A635").NumberFormat = "@"
$ws.Range("A635").Value = "2026/01/16"
$ws.Range("A635").ClearFormats()

$ws.Range("B635").Value = "金"
$ws.Range("C635").Value = 14
$ws.Range("D635").Value = 58
